$wb = $excel.ActiveWorkbook

# --- New spider monster entries on the "Monsters" sheet ---
$ws = $wb.Worksheets.Item("Monsters")

# Single monster row (columns A-C): id 58 "Höhlenspinne"
$ws.Cells.Item(3, 1).Value = 58
$ws.Cells.Item(3, 2).Value = "Höhlenspinne"
$ws.Cells.Item(3, 3).Value = "Monster in Ship's end"

# Monster group rows (columns H-J): group ids 88 and 89
$ws.Cells.Item(3, 8).Value = 88
$ws.Cells.Item(3, 9).Value = "3x Höhlenspinne"

$ws.Cells.Item(4, 8).Value = 89
$ws.Cells.Item(4, 9).Value = "4x Höhlenspinne"

# --- Switch the active/selected sheet from "Todo" to "Monsters" ---
$ws.Activate() | Out-Null
$ws.Range("J4").Select() | Out-Null
